# update master-pegwai and login
# The "password" column (column B) is no longer part of the exported
# pegawai/login master data, so remove it entirely - this shifts every
# column after it one position to the left and drops the trailing
# "analissdm" column that fell off the end (N1 -> gone, now ends at M1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column B first (mirrors how a user would click the column header
# before deleting it), then delete the whole column so everything to the
# right shifts left.
$ws.Columns.Item(2).Select() | Out-Null
$ws.Columns.Item(2).Delete() | Out-Null
